$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Swap the three tables (slides 14, 15, 16) from the custom table style
#    {5BBFCD71-6738-42E8-9DEF-8367C8461830} to the built-in style
#    {91230868-9EC4-46BF-A3CC-6F292B9285E9}. In each of these slides the
#    table is the first shape on the slide.
# ---------------------------------------------------------------------------
$newTableStyle = "{91230868-9EC4-46BF-A3CC-6F292B9285E9}"
foreach ($slideIdx in 14,15,16) {
    $slide = $p.Slides.Item($slideIdx)
    $shp = $slide.Shapes.Item(1)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle($newTableStyle)
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the presentation's (slide master) theme from the "Integral /
#    Red Violet" palette to the stock "Office" palette (the 12 scheme
#    colours are the only part of the theme that actually differs between
#    the two themes - fonts and format scheme are already identical).
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$tcs = $slide1.ThemeColorScheme
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
